# Fix "Unclaimed Property" Year values on the NewTaxReturn sheet so they are
# consistent with the rest of their block (2024 / 2023 / 2022), and add the
# missing CRN ("Y") flag for the tax types that were missing it in each
# year block: Motor Fuel Tax, Tire Recycling Fee, Tobacco Tax and
# Transportation Network Services.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTaxReturn")

# --- Correct the Year ("Month"/date) values for Unclaimed Property rows ---
$ws.Range("F27").Value = "2024"
$ws.Range("F39").Value = "2023"
$ws.Range("F51").Value = "2022"

# --- Add the missing CRN flags ---
$crnRows = @(21, 24, 25, 26, 33, 36, 37, 38, 45, 48, 49, 50)
foreach ($r in $crnRows) {
    $ws.Range("H$r").Value = "Y"
}
